{"js": "// Replace the date line and every \"A\u00f7B=C, D\" division answer in the table\n// with the values from the updated answer key. Each old string occurs\n// exactly once in the document, so a scoped search-and-replace per pair is\n// safe and avoids having to hard-code table/row/column coordinates.\nconst replacements = [\n  [\"2025-05-29 Thursday\", \"2025-05-30 Friday\"],\n  [\"221\u00f75=44, 1\", \"735\u00f74=183, 3\"],\n  [\"176\u00f72=88, 0\", \"675\u00f74=168, 3\"],\n  [\"490\u00f79=54, 4\", \"955\u00f75=191, 0\"],\n  [\"665\u00f79=73, 8\", \"963\u00f75=192, 3\"],\n  [\"768\u00f74=192, 0\", \"354\u00f79=39, 3\"],\n  [\"728\u00f76=121, 2\", \"722\u00f77=103, 1\"],\n  [\"547\u00f76=91, 1\", \"615\u00f79=68, 3\"],\n  [\"608\u00f77=86, 6\", \"352\u00f77=50, 2\"],\n  [\"937\u00f77=133, 6\", \"643\u00f77=91, 6\"],\n  [\"692\u00f79=76, 8\", \"855\u00f76=142, 3\"],\n  [\"783\u00f77=111, 6\", \"256\u00f74=64, 0\"],\n  [\"754\u00f74=188, 2\", \"914\u00f74=228, 2\"],\n  [\"230\u00f78=28, 6\", \"341\u00f79=37, 8\"],\n  [\"407\u00f79=45, 2\", \"504\u00f74=126, 0\"],\n  [\"186\u00f72=93, 0\", \"133\u00f73=44, 1\"],\n  [\"863\u00f75=172, 3\", \"639\u00f75=127, 4\"],\n  [\"939\u00f73=313, 0\", \"523\u00f78=65, 3\"],\n  [\"892\u00f78=111, 4\", \"731\u00f74=182, 3\"],\n  [\"982\u00f75=196, 2\", \"368\u00f78=46, 0\"],\n  [\"228\u00f74=57, 0\", \"211\u00f76=35, 1\"],\n  [\"209\u00f72=104, 1\", \"323\u00f73=107, 2\"],\n  [\"259\u00f72=129, 1\", \"329\u00f72=164, 1\"],\n  [\"456\u00f75=91, 1\", \"435\u00f79=48, 3\"],\n  [\"532\u00f74=133, 0\", \"940\u00f77=134, 2\"],\n  [\"256\u00f77=36, 4\", \"284\u00f73=94, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of found.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00f7B=C, D\" division answer in the table\n# with the values from the updated answer key. Each old string occurs\n# exactly once in the document, so a Find/Replace pass per pair is safe and\n# avoids having to hard-code table/row/column coordinates.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2025-05-29 Thursday\"; Replace = \"2025-05-30 Friday\" },\n    @{ Find = \"221\u00f75=44, 1\"; Replace = \"735\u00f74=183, 3\" },\n    @{ Find = \"176\u00f72=88, 0\"; Replace = \"675\u00f74=168, 3\" },\n    @{ Find = \"490\u00f79=54, 4\"; Replace = \"955\u00f75=191, 0\" },\n    @{ Find = \"665\u00f79=73, 8\"; Replace = \"963\u00f75=192, 3\" },\n    @{ Find = \"768\u00f74=192, 0\"; Replace = \"354\u00f79=39, 3\" },\n    @{ Find = \"728\u00f76=121, 2\"; Replace = \"722\u00f77=103, 1\" },\n    @{ Find = \"547\u00f76=91, 1\"; Replace = \"615\u00f79=68, 3\" },\n    @{ Find = \"608\u00f77=86, 6\"; Replace = \"352\u00f77=50, 2\" },\n    @{ Find = \"937\u00f77=133, 6\"; Replace = \"643\u00f77=91, 6\" },\n    @{ Find = \"692\u00f79=76, 8\"; Replace = \"855\u00f76=142, 3\" },\n    @{ Find = \"783\u00f77=111, 6\"; Replace = \"256\u00f74=64, 0\" },\n    @{ Find = \"754\u00f74=188, 2\"; Replace = \"914\u00f74=228, 2\" },\n    @{ Find = \"230\u00f78=28, 6\"; Replace = \"341\u00f79=37, 8\" },\n    @{ Find = \"407\u00f79=45, 2\"; Replace = \"504\u00f74=126, 0\" },\n    @{ Find = \"186\u00f72=93, 0\"; Replace = \"133\u00f73=44, 1\" },\n    @{ Find = \"863\u00f75=172, 3\"; Replace = \"639\u00f75=127, 4\" },\n    @{ Find = \"939\u00f73=313, 0\"; Replace = \"523\u00f78=65, 3\" },\n    @{ Find = \"892\u00f78=111, 4\"; Replace = \"731\u00f74=182, 3\" },\n    @{ Find = \"982\u00f75=196, 2\"; Replace = \"368\u00f78=46, 0\" },\n    @{ Find = \"228\u00f74=57, 0\"; Replace = \"211\u00f76=35, 1\" },\n    @{ Find = \"209\u00f72=104, 1\"; Replace = \"323\u00f73=107, 2\" },\n    @{ Find = \"259\u00f72=129, 1\"; Replace = \"329\u00f72=164, 1\" },\n    @{ Find = \"456\u00f75=91, 1\"; Replace = \"435\u00f79=48, 3\" },\n    @{ Find = \"532\u00f74=133, 0\"; Replace = \"940\u00f77=134, 2\" },\n    @{ Find = \"256\u00f77=36, 4\"; Replace = \"284\u00f73=94, 2\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.Find\n    $found = $find.Execute(\n        $r.Find,    # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $r.Replace, # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        Write-Output \"WARNING: not found: $($r.Find)\"\n    }\n}\n\nWrite-Output \"done\"\n"}
